# Apply "Add data for 2022-10-12" update to the carjacking-by-neighborhood-by-month
# workbook: rename the sheet / update the header label to reflect the new
# "through" date, and update/insert the affected monthly counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab and update the "through" date text in the header row.
$ws.Name = "Through 2022-10-04"
$ws.Range("B1").Value = "October 2022 (through October 04)"

# --- Garfield Park (row 2) ---
$ws.Range("V2").Value = 3

# --- Humboldt Park (row 3) ---
$ws.Range("L3").Value = 2

# --- Englewood (row 4) ---
$ws.Range("V4").Value = 2
$ws.Range("AF4").Value = 1
$ws.Range("AZ4").Value = 2

# --- Austin (row 5) ---
$ws.Range("V5").Value = 2
$ws.Range("AP5").Value = 2

# --- North Lawndale (row 6) ---
$ws.Range("L6").Value = 3
$ws.Range("V6").Value = 3
$ws.Range("BT6").Value = 1

# --- South Shore (row 7) ---
$ws.Range("AP7").Value = 1

# --- Roseland (row 13) ---
$ws.Range("B13").Value = 1
$ws.Range("V13").Value = 1

# --- Chatham (row 16) ---
$ws.Range("B16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("AF16").Value = 1

# --- Washington Heights (row 17) ---
$ws.Range("V17").Value = 2

# --- Grand Boulevard (row 19) ---
$ws.Range("L19").Value = 1

# --- Lower West Side (row 21) ---
$ws.Range("L21").Value = 2

# --- Ashburn (row 22) ---
$ws.Range("V22").Value = 1

# --- Auburn Gresham (row 23) ---
$ws.Range("L23").Value = 4

# --- Logan Square (row 27) ---
$ws.Range("V27").Value = 2

# --- Lake View (row 29) ---
$ws.Range("L29").Value = 1
$ws.Range("AP29").Value = 1

# --- West Town (row 33) ---
$ws.Range("L33").Value = 3

# --- Bucktown (row 47) ---
$ws.Range("B47").Value = 1

# --- Calumet Heights (row 48) ---
$ws.Range("L48").Value = 1

# --- Chicago Lawn (row 66) ---
$ws.Range("V66").Value = 1

# --- Hermosa (row 78) ---
$ws.Range("AP78").Value = 2

# --- Lincoln Square (row 82) ---
$ws.Range("L82").Value = 1

# --- North Center (row 87) ---
$ws.Range("L87").Value = 1
$ws.Range("AP87").Value = 1

# --- Norwood Park (row 88) ---
$ws.Range("L88").Value = 1
